# Atualização de bases das ligas, do dia: 11-04-2024 às 00:31
#
# For each listed row-pair, the data (columns B through AC) got swapped
# between the two rows while the index column (A) stayed the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

$rowPairs = @(
    @(14, 15),
    @(130, 131),
    @(138, 139)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"

        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2

        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}
